$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('B12').Value = 'Type 02: demand in imperative form + threat in active form'
$ws.Range('B23').Value = 'Type 02: demand in imperative form + threat in active form'
$ws.Range('B57').Value = 'Type 03: threat in active form + demand in active form'
$ws.Range('B95').Value = 'Type 03: threat in active form + demand in active form'
$ws.Range('B111').Value = 'Type 03: threat in active form + demand in active form'
$ws.Range('B113').Value = 'Type 02: demand in imperative form + threat in active form'
$ws.Range('B159').Value = 'Type 03: threat in active form + demand in active form'
$ws.Range('B161').Value = 'Type 03: threat in active form + demand in active form'
$ws.Range('B162').Value = 'Type 03: threat in active form + demand in active form'
$ws.Range('B164').Value = 'Type 03: threat in active form + demand in active form'
$ws.Range('B181').Value = 'Type 03: threat in active form + demand in active form'
$ws.Range('B183').Value = 'Type 03: threat in active form + demand in active form'
$ws.Range('B184').Value = 'Type 03: threat in active form + demand in active form'
$ws.Range('B187').Value = 'Type 03: threat in active form + demand in active form'
$ws.Range('B195').Value = 'Type 03: threat in active form + demand in active form'
$ws.Range('B197').Value = 'Type 03: threat in active form + demand in active form'
$ws.Range('B198').Value = 'Type 03: threat in active form + demand in active form'
$ws.Range('B200').Value = 'Type 03: threat in active form + demand in active form'
$ws.Range('B204').Value = 'Type 03: threat in active form + demand in active form'
$ws.Range('B213').Value = 'Type 03: threat in active form + demand in active form'
$ws.Range('B214').Value = 'Type 03: threat in active form + demand in active form'
$ws.Range('B218').Value = 'Type 03: threat in active form + demand in active form'
$ws.Range('B220').Value = 'Type 03: threat in active form + demand in active form'
$ws.Range('B221').Value = 'Type 03: threat in active form + demand in active form'
$ws.Range('B227').Value = 'Type 03: threat in active form + demand in active form'
$ws.Range('B229').Value = 'Type 03: threat in active form + demand in active form'
$ws.Range('B230').Value = 'Type 03: threat in active form + demand in active form'
$ws.Range('B231').Value = 'Type 03: threat in active form + demand in active form'
$ws.Range('B324').Value = 'Type 03: threat in active form + demand in active form'
$ws.Range('B366').Value = 'Type 02: demand in imperative form + threat in active form'
$ws.Range('B457').Value = 'Type 03: threat in active form + demand in active form'
$ws.Range('B471').Value = 'Type 03: threat in active form + demand in active form'
$ws.Range('B569').Value = 'Type 03: threat in active form + demand in active form'
$ws.Range('B740').Value = 'Type 03: threat in active form + demand in active form'
$ws.Range('B750').Value = 'Type 03: threat in active form + demand in active form'
$ws.Range('B756').Value = 'No Extortion'
$ws.Range('B765').Value = 'No Extortion'
$ws.Range('B774').Value = 'Type 03: threat in active form + demand in active form'
$ws.Range('B784').Value = 'No Extortion'
$ws.Range('B796').Value = 'No Extortion'
$ws.Range('B800').Value = 'Type 03: threat in active form + demand in active form'
$ws.Range('B806').Value = 'Type 03: threat in active form + demand in active form'
$ws.Range('B807').Value = 'Type 03: threat in active form + demand in active form'
$ws.Range('B812').Value = 'Type 03: threat in active form + demand in active form'
$ws.Range('B816').Value = 'No Extortion'
$ws.Range('B819').Value = 'No Extortion'
$ws.Range('B822').Value = 'No Extortion'
$ws.Range('B844').Value = 'Type 03: threat in active form + demand in active form'
$ws.Range('B847').Value = 'Type 03: threat in active form + demand in active form'
$ws.Range('B863').Value = 'No Extortion'
$ws.Range('B871').Value = 'Type 03: threat in active form + demand in active form'
$ws.Range('B880').Value = 'Type 03: threat in active form + demand in active form'
$ws.Range('B894').Value = 'Type 03: threat in active form + demand in active form'
$ws.Range('B895').Value = 'No Extortion'
$ws.Range('B908').Value = 'No Extortion'
$ws.Range('B919').Value = 'No Extortion'
$ws.Range('B925').Value = 'No Extortion'
$ws.Range('B929').Value = 'Type 03: threat in active form + demand in active form'
$ws.Range('B947').Value = 'Type 03: threat in active form + demand in active form'
$ws.Range('B950').Value = 'Type 03: threat in active form + demand in active form'
$ws.Range('B953').Value = 'Type 03: threat in active form + demand in active form'
$ws.Range('B973').Value = 'Type 03: threat in active form + demand in active form'
$ws.Range('B974').Value = 'Type 03: threat in active form + demand in active form'
$ws.Range('B975').Value = 'Type 03: threat in active form + demand in active form'
$ws.Range('B982').Value = 'Type 03: threat in active form + demand in active form'
